$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("data_info")
$sheet2 = $wb.Worksheets.Item("timecards_log")

# --- data_info sheet: clarify the description for the work_date column ---
$sheet1.Range("D4").Value = "Work date for the timecard entry (YYYY-MM-DD)"

# --- timecards_log sheet: expand the single summary row into a full year of
#     monthly rows, and update the findings / change-log notes ---

# Monthly summary data (month start date serial, rows, missing start, missing
# end, missing lunch start, missing lunch end, change-log note)
$months = @(
    @{ Date = 45292; Rows = 400; Start = 0; End = 0; LunchStart = 103; LunchEnd = 103; Note = $null },
    @{ Date = 45323; Rows = 400; Start = 0; End = 0; LunchStart = 118; LunchEnd = 118; Note = $null },
    @{ Date = 45352; Rows = 400; Start = 0; End = 0; LunchStart = 112; LunchEnd = 112; Note = $null },
    @{ Date = 45383; Rows = 400; Start = 0; End = 0; LunchStart = 125; LunchEnd = 125; Note = $null },
    @{ Date = 45413; Rows = 400; Start = 0; End = 0; LunchStart = 108; LunchEnd = 108; Note = $null },
    @{ Date = 45444; Rows = 400; Start = 7; End = 6; LunchStart = 190; LunchEnd = 191; Note = "13 rows removed for verificaiton" },
    @{ Date = 45474; Rows = 460; Start = 9; End = 9; LunchStart = 170; LunchEnd = 166; Note = "17 rows removed for verificaiton" },
    @{ Date = 45505; Rows = 400; Start = 0; End = 0; LunchStart = 119; LunchEnd = 119; Note = $null },
    @{ Date = 45536; Rows = 400; Start = 0; End = 0; LunchStart = 121; LunchEnd = 121; Note = $null },
    @{ Date = 45566; Rows = 400; Start = 0; End = 0; LunchStart = 112; LunchEnd = 112; Note = $null },
    @{ Date = 45597; Rows = 400; Start = 0; End = 0; LunchStart = 119; LunchEnd = 119; Note = $null },
    @{ Date = 45627; Rows = 400; Start = 0; End = 0; LunchStart = 131; LunchEnd = 131; Note = $null }
)

$r = 2
foreach ($m in $months) {
    $sheet2.Cells.Item($r, 1).Value = $m.Date
    $sheet2.Cells.Item($r, 1).HorizontalAlignment = -4108
    $sheet2.Cells.Item($r, 1).VerticalAlignment = -4160

    $sheet2.Cells.Item($r, 2).Value = $m.Rows
    $sheet2.Cells.Item($r, 3).Value = "yes"
    $sheet2.Cells.Item($r, 4).Value = $m.Start
    $sheet2.Cells.Item($r, 5).Value = $m.End
    $sheet2.Cells.Item($r, 6).Value = $m.LunchStart
    $sheet2.Cells.Item($r, 7).Value = $m.LunchEnd

    if ($m.Note) {
        $sheet2.Cells.Item($r, 8).Value = $m.Note
    } else {
        $sheet2.Cells.Item($r, 8).Value = ""
    }
    $sheet2.Cells.Item($r, 8).HorizontalAlignment = -4152
    $sheet2.Cells.Item($r, 8).VerticalAlignment = -4160

    $r = $r + 1
}

# --- switch the active tab to the timecards_log sheet ---
$sheet2.Activate()
$sheet2.Range("D22").Select()
